# Updated cryptos list values per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Preserve literal text (e.g. "79.80", "1.00") instead of letting
    # Excel auto-coerce numeric-looking strings into numbers, which would
    # drop trailing zeros / thousands-dot formatting. Restore the original
    # cell style afterwards so no stray formatting is introduced.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "69.815.41"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.680.95"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "650.44"
$ws.Range("E5").Value = "  -4.10%  "
Set-TextValue $ws.Range("D6") "161.38"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.14%  "
Set-TextValue $ws.Range("D8") "0.498"
$ws.Range("E8").Value = "  +0.35%  "
Set-TextValue $ws.Range("D9") "0.145"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  +0.27%  "
Set-TextValue $ws.Range("D11") "0.442"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").Value = "4.309.39"
$ws.Range("E13").Value = "  -0.49%  "
Set-TextValue $ws.Range("D14") "32.69"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "3.656.97"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "69.807.29"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  +0.34%  "
Set-TextValue $ws.Range("D18") "15.95"
$ws.Range("E18").Value = "  -0.81%  "
Set-TextValue $ws.Range("D19") "6.52"
$ws.Range("E19").Value = "  +0.18%  "
Set-TextValue $ws.Range("D20") "10.34"
$ws.Range("E20").Value = "  +5.44%  "
Set-TextValue $ws.Range("D21") "471.03"
$ws.Range("E21").Value = "  -0.52%  "
Set-TextValue $ws.Range("D22") "0.652"
$ws.Range("E22").Value = "  -0.09%  "
Set-TextValue $ws.Range("D23") "79.80"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "3.833.04"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E25").Value = "  +0.22%  "
Set-TextValue $ws.Range("D27") "11.03"
$ws.Range("E27").Value = "  +0.45%  "
Set-TextValue $ws.Range("D28") "8.83"
$ws.Range("E28").Value = "  -3.30%  "
Set-TextValue $ws.Range("D29") "2.65"
$ws.Range("E29").Value = "  -2.20%  "
Set-TextValue $ws.Range("D30") "1.70"
$ws.Range("E30").Value = "  -2.80%  "
Set-TextValue $ws.Range("D31") "2.00"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D32") "1.00"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D33") "0.166"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D34") "6.52"
$ws.Range("E34").Value = "  -1.23%  "
Set-TextValue $ws.Range("D35") "26.73"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").Value = "3.679.59"
$ws.Range("E36").Value = "  -0.44%  "
Set-TextValue $ws.Range("D37") "8.40"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("E38").Value = "  -0.09%  "
Set-TextValue $ws.Range("D39") "5.88"
$ws.Range("E39").Value = "  -5.37%  "
Set-TextValue $ws.Range("D40") "178.40"
$ws.Range("E40").Value = "  +6.92%  "
Set-TextValue $ws.Range("D42") "2.21"
$ws.Range("E42").Value = "  -1.55%  "
Set-TextValue $ws.Range("D43") "0.0898"
$ws.Range("E43").Value = "  -0.81%  "
Set-TextValue $ws.Range("D44") "0.929"
$ws.Range("E44").Value = "  -1.77%  "
Set-TextValue $ws.Range("D45") "2.83"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D46") "46.53"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "28.95"
$ws.Range("E47").Value = "  +1.68%  "
Set-TextValue $ws.Range("D48") "0.000272"
$ws.Range("E48").Value = "  -3.11%  "
Set-TextValue $ws.Range("D49") "7.86"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("E51").Value = "  -5.98%  "
